$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update image file names in column A (rows 2-4)
$ws.Range("A2").Value = "im1.jpg"
$ws.Range("A3").Value = "im2.jpg"
$ws.Range("A4").Value = "im3.jpg"

# Move the active selection to A5
$ws.Range("A5").Select()
